$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (data no longer needed)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "F13a1"
$ws.Range("C2").Value = "Itga9"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07623033333333334
$ws.Range("H2").Value = 0.228691
$ws.Range("I2").Value = 0.7411556909515168
$ws.Range("J2").Value = 0.7411556909515167
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7521946666666667
$ws.Range("N2").Value = 2.256584
$ws.Range("O2").Value = 0.07361670343069449
$ws.Range("P2").Value = 0.0736167034306945
$ws.Range("Q2").Value = 0.05734005017155556
$ws.Range("R2").Value = 0.5160604515440002
$ws.Range("S2").Value = 0.05456143869674927
$ws.Range("T2").Value = 0.05456143869674927

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "F13a1"
$ws.Range("C3").Value = "Itga9"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07623033333333334
$ws.Range("H3").Value = 0.228691
$ws.Range("I3").Value = 0.7411556909515168
$ws.Range("J3").Value = 0.7411556909515167
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.059280333333334
$ws.Range("N3").Value = 21.177841
$ws.Range("O3").Value = 0.6908862423022597
$ws.Range("P3").Value = 0.6908862423022598
$ws.Range("Q3").Value = 0.5381312929034445
$ws.Range("R3").Value = 4.843181636131001
$ws.Range("S3").Value = 0.5120542702824283
$ws.Range("T3").Value = 0.5120542702824283

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F13a1"
$ws.Range("C4").Value = "Itga9"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07623033333333334
$ws.Range("H4").Value = 0.228691
$ws.Range("I4").Value = 0.7411556909515168
$ws.Range("J4").Value = 0.7411556909515167
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.406242333333334
$ws.Range("N4").Value = 7.218727
$ws.Range("O4").Value = 0.2354970542670457
$ws.Range("P4").Value = 0.2354970542670457
$ws.Range("Q4").Value = 0.1834286551507778
$ws.Range("R4").Value = 1.650857896357
$ws.Range("S4").Value = 0.1745399819723391
$ws.Range("T4").Value = 0.1745399819723391

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "F13a1"
$ws.Range("C5").Value = "Itga9"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.026623
$ws.Range("H5").Value = 0.079869
$ws.Range("I5").Value = 0.2588443090484832
$ws.Range("J5").Value = 0.2588443090484832
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7521946666666667
$ws.Range("N5").Value = 2.256584
$ws.Range("O5").Value = 0.07361670343069449
$ws.Range("P5").Value = 0.0736167034306945
$ws.Range("Q5").Value = 0.02002567861066666
$ws.Range("R5").Value = 0.180231107496
$ws.Range("S5").Value = 0.01905526473394522
$ws.Range("T5").Value = 0.01905526473394522

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F13a1"
$ws.Range("C6").Value = "Itga9"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.026623
$ws.Range("H6").Value = 0.079869
$ws.Range("I6").Value = 0.2588443090484832
$ws.Range("J6").Value = 0.2588443090484832
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.059280333333334
$ws.Range("N6").Value = 21.177841
$ws.Range("O6").Value = 0.6908862423022597
$ws.Range("P6").Value = 0.6908862423022598
$ws.Range("Q6").Value = 0.1879392203143333
$ws.Range("R6").Value = 1.691452982829
$ws.Range("S6").Value = 0.1788319720198314
$ws.Range("T6").Value = 0.1788319720198314

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F13a1"
$ws.Range("C7").Value = "Itga9"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.026623
$ws.Range("H7").Value = 0.079869
$ws.Range("I7").Value = 0.2588443090484832
$ws.Range("J7").Value = 0.2588443090484832
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.406242333333334
$ws.Range("N7").Value = 7.218727
$ws.Range("O7").Value = 0.2354970542670457
$ws.Range("P7").Value = 0.2354970542670457
$ws.Range("Q7").Value = 0.06406138964033334
$ws.Range("R7").Value = 0.576552506763
$ws.Range("S7").Value = 0.06095707229470659
$ws.Range("T7").Value = 0.06095707229470659
